# "Export excel overall sale actual."
# Clears the actual-sales figures on the "1.Overall Monthly sales result"
# sheet (the source data behind the "Over All Sales Result" chart) back to
# zero, e.g. before handing the template out / exporting a blank report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1.Overall Monthly sales result")

# Rows: 3 = OTHER, 4 = PMSP, 5 = OEM, 7 = 19OAP. Row 6 (TOTAL) is a SUM()
# formula over B3:B5 etc. and recalculates to 0 on its own.
$ws.Range("B3:E5").Value = 0
$ws.Range("B7:E7").Value = 0

# Leave the selection where the author left it when they saved.
$ws.Range("P12").Select()

$wb.Save()
